$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New set of words that replaces the previous "working set" of German verbs
# used in the pairings sequence (column C, rows 2-129).
$newWords = @(
    'zeigen',
    'reizen',
    'formen',
    'decken',
    'folgen',
    'tollen',
    'heulen',
    'pflanzen',
    'spinnen',
    'achten',
    'boxen',
    'saufen',
    'freuen',
    'zünden',
    'wachsen',
    'dringen',
    'drehen',
    'jubeln',
    'jagen',
    'fallen',
    'fischen',
    'scheinen',
    'lockern',
    'gründen',
    'spüren',
    'wüten',
    'führen',
    'bauen',
    'sterben',
    'malen',
    'hauen',
    'fahren',
    'mögen',
    'pfeifen',
    'wenden',
    'töten',
    'flüchten',
    'zögern',
    'erben',
    'biegen',
    'räumen',
    'stehlen',
    'stammen',
    'schwingen',
    'fangen',
    'seufzen',
    'loben',
    'schlucken',
    'werfen',
    'kichern',
    'irren',
    'zielen',
    'dienen',
    'scheitern',
    'lügen',
    'gelten',
    'ärgern',
    'sprengen',
    'bergen',
    'sichern',
    'weichen',
    'streichen',
    'äußern',
    'quälen',
    'wirken',
    'heben',
    'heilen',
    'siegen',
    'bellen',
    'kümmern',
    'filmen',
    'wehtun',
    'hören',
    'stecken',
    'kosten',
    'betteln',
    'rufen',
    'werden',
    'runden',
    'tropfen',
    'schreiten',
    'treiben',
    'münzen',
    'fällen',
    'klettern',
    'planen',
    'knarren',
    'trennen',
    'liegen',
    'sperren',
    'machen',
    'greifen',
    'fließen',
    'backen',
    'schwören',
    'kehren',
    'trauen',
    'bitten',
    'sinken',
    'warnen',
    'suchen',
    'schrecken',
    'graben',
    'rasen',
    'liefern',
    'ändern',
    'wundern',
    'lesen',
    'spielen',
    'arten',
    'klingen',
    'schulden',
    'brauchen',
    'enden',
    'grüßen',
    'schwächen',
    'geben',
    'bluten',
    'platzen',
    'mauern',
    'sorgen',
    'schenken',
    'schmecken',
    'ehren',
    'altern',
    'feiern',
    'helfen',
    'flehen'
)

$startRow = 2
for ($i = 0; $i -lt $newWords.Count; $i++) {
    $ws.Cells.Item($startRow + $i, 3).Value = $newWords[$i]
}
